$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 606028.7
$ws.Range("I15").Value = 606028.7
$ws.Range("K15").Value = 1818086.1
$ws.Range("M15").Value = -1817917.1
$ws.Range("H18").Value = 8049.727
$ws.Range("I18").Value = 7588.5557
$ws.Range("K18").Value = 7588.5557
$ws.Range("M18").Value = -7304.5557
$ws.Range("H40").Value = 5600
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H96").Value = 826.4
$ws.Range("I96").Value = 600
$ws.Range("J96").Value = 923.4286
$ws.Range("K96").Value = 1800
$ws.Range("L96").Value = 2770.2858
$ws.Range("M96").Value = -427
$ws.Range("N96").Value = -5516.2858
$ws.Range("H98").Value = 54599.273
$ws.Range("I98").Value = 69449.625
$ws.Range("K98").Value = 69449.625
$ws.Range("M98").Value = -67951.625
$ws.Range("H112").Value = 2681
$ws.Range("J112").Value = 1501.6666
$ws.Range("L112").Value = 4504.9998
$ws.Range("N112").Value = -6720.9998
$ws.Range("H122").Value = 54599.273
$ws.Range("I122").Value = 69449.625
$ws.Range("K122").Value = 208348.875
$ws.Range("M122").Value = -205898.875
$ws.Range("H132").Value = 1889987.5
$ws.Range("I132").Value = 2986.8
$ws.Range("K132").Value = 8960.400000000001
$ws.Range("M132").Value = -6430.400000000001
$ws.Range("H137").Value = 639960.5600000001
$ws.Range("I137").Value = 824621.9399999999
$ws.Range("J137").Value = 3904.6667
$ws.Range("K137").Value = 2473865.82
$ws.Range("L137").Value = 11714.0001
$ws.Range("M137").Value = -2471315.82
$ws.Range("N137").Value = -16814.0001
$ws.Range("H138").Value = 166825.52
$ws.Range("I138").Value = 628921.5
$ws.Range("J138").Value = 5091.9165
$ws.Range("K138").Value = 1886764.5
$ws.Range("L138").Value = 15275.7495
$ws.Range("M138").Value = -1881624.5
$ws.Range("N138").Value = -25555.7495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11290.366
$ws.Range("I32").Value = 10307.459
$ws.Range("J32").Value = 27999.8
$ws.Range("K32").Value = 10307.459
$ws.Range("L32").Value = 27999.8
$ws.Range("M32").Value = -10020.459
$ws.Range("N32").Value = -28573.8
$ws.Range("H63").Value = 3666.3333
$ws.Range("I63").Value = 3999.5
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 3999.5
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -3313.5
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 3666.3333
$ws.Range("I66").Value = 3999.5
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 19997.5
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -16565.5
$ws.Range("N66").Value = -21864
$ws.Range("H74").Value = 3896.4482
$ws.Range("I74").Value = 20949
$ws.Range("K74").Value = 20949
$ws.Range("M74").Value = -20075
$ws.Range("H77").Value = 3896.4482
$ws.Range("I77").Value = 20949
$ws.Range("K77").Value = 104745
$ws.Range("M77").Value = -100377
$ws.Range("H122").Value = 1006005.75
$ws.Range("I122").Value = 5968.8335
$ws.Range("J122").Value = 5006153.5
$ws.Range("K122").Value = 17906.5005
$ws.Range("L122").Value = 15018460.5
$ws.Range("M122").Value = -15456.5005
$ws.Range("N122").Value = -15023360.5
$ws.Range("H132").Value = 1646.65
$ws.Range("I132").Value = 1313.3143
$ws.Range("K132").Value = 3939.9429
$ws.Range("M132").Value = -1409.9429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19970.9
$ws.Range("I82").Value = 9963.625
$ws.Range("J82").Value = 60000
$ws.Range("K82").Value = 9963.625
$ws.Range("L82").Value = 60000
$ws.Range("M82").Value = -9580.625
$ws.Range("N82").Value = -60766
$ws.Range("H85").Value = 19970.9
$ws.Range("I85").Value = 9963.625
$ws.Range("J85").Value = 60000
$ws.Range("K85").Value = 9963.625
$ws.Range("L85").Value = 60000
$ws.Range("M85").Value = -8637.625
$ws.Range("N85").Value = -62652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2324.375
$ws.Range("I16").Value = 1219
$ws.Range("J16").Value = 4166.6665
$ws.Range("K16").Value = 1219
$ws.Range("L16").Value = 4166.6665
$ws.Range("M16").Value = -932
$ws.Range("N16").Value = -4740.6665
$ws.Range("H22").Value = 381.16666
$ws.Range("I22").Value = 419
$ws.Range("J22").Value = 305.5
$ws.Range("K22").Value = 419
$ws.Range("L22").Value = 305.5
$ws.Range("M22").Value = -69
$ws.Range("N22").Value = -1005.5
$ws.Range("H31").Value = 4111.1665
$ws.Range("I31").Value = 1890.8334
$ws.Range("J31").Value = 6331.5
$ws.Range("K31").Value = 1890.8334
$ws.Range("L31").Value = 6331.5
$ws.Range("M31").Value = -1595.8334
$ws.Range("N31").Value = -6921.5
$ws.Range("H34").Value = 4111.1665
$ws.Range("I34").Value = 1890.8334
$ws.Range("J34").Value = 6331.5
$ws.Range("K34").Value = 1890.8334
$ws.Range("L34").Value = 6331.5
$ws.Range("M34").Value = -1688.8334
$ws.Range("N34").Value = -6735.5
$ws.Range("H41").Value = 10000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H113").Value = 2324.375
$ws.Range("I113").Value = 1219
$ws.Range("J113").Value = 4166.6665
$ws.Range("K113").Value = 1219
$ws.Range("L113").Value = 4166.6665
$ws.Range("M113").Value = 951
$ws.Range("N113").Value = -8506.666499999999
$ws.Range("H141").Value = 457876.84
$ws.Range("J141").Value = 613890.75
$ws.Range("L141").Value = 613890.75
$ws.Range("N141").Value = -624250.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1958
$ws.Range("I63").Value = 1266.6666
$ws.Range("J63").Value = 2995
$ws.Range("K63").Value = 3799.9998
$ws.Range("L63").Value = 8985
$ws.Range("M63").Value = -3050.9998
$ws.Range("N63").Value = -10483
$ws.Range("H66").Value = 1958
$ws.Range("I66").Value = 1266.6666
$ws.Range("J66").Value = 2995
$ws.Range("K66").Value = 11399.9994
$ws.Range("L66").Value = 26955
$ws.Range("M66").Value = -7655.999400000001
$ws.Range("N66").Value = -34443
$ws.Range("H74").Value = 17066.334
$ws.Range("I74").Value = 4997.5
$ws.Range("J74").Value = 18923.076
$ws.Range("K74").Value = 14992.5
$ws.Range("L74").Value = 56769.228
$ws.Range("M74").Value = -13931.5
$ws.Range("N74").Value = -58891.228
$ws.Range("H75").Value = 869.125
$ws.Range("I75").Value = 842.75
$ws.Range("K75").Value = 2528.25
$ws.Range("M75").Value = -1530.25
$ws.Range("H77").Value = 17066.334
$ws.Range("I77").Value = 4997.5
$ws.Range("J77").Value = 18923.076
$ws.Range("K77").Value = 44977.5
$ws.Range("L77").Value = 170307.684
$ws.Range("M77").Value = -39673.5
$ws.Range("N77").Value = -180915.684
$ws.Range("H78").Value = 869.125
$ws.Range("I78").Value = 842.75
$ws.Range("K78").Value = 7584.75
$ws.Range("M78").Value = -2592.75
$ws.Range("H132").Value = 49961.383
$ws.Range("J132").Value = 53974.832
$ws.Range("L132").Value = 485773.488
$ws.Range("N132").Value = -490833.488

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 14183.9375
$ws.Range("I122").Value = 10149.462
$ws.Range("J122").Value = 31666.666
$ws.Range("K122").Value = 30448.386
$ws.Range("L122").Value = 94999.99800000001
$ws.Range("M122").Value = -27998.386
$ws.Range("N122").Value = -99899.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1205.2941
$ws.Range("I46").Value = 959.9
$ws.Range("K46").Value = 959.9
$ws.Range("M46").Value = -771.9
$ws.Range("H64").Value = 16157.833
$ws.Range("J64").Value = 16157.833
$ws.Range("L64").Value = 16157.833
$ws.Range("N64").Value = -16607.833
$ws.Range("H67").Value = 16157.833
$ws.Range("J67").Value = 16157.833
$ws.Range("L67").Value = 16157.833
$ws.Range("N67").Value = -17717.833
$ws.Range("H138").Value = 64393.75
$ws.Range("I138").Value = 64388
$ws.Range("J138").Value = 64395.668
$ws.Range("K138").Value = 64388
$ws.Range("L138").Value = 64395.668
$ws.Range("M138").Value = -59248
$ws.Range("N138").Value = -74675.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5980.0347
$ws.Range("I122").Value = 4481.1816
$ws.Range("K122").Value = 13443.5448
$ws.Range("M122").Value = -10993.5448
